$d = $word.ActiveDocument

# --- Change 1: merge the "Christian " / "Ehringfeld" / ", David Mock, ..." runs
# (and drop the spell-check proofErr markers around "Ehringfeld") into a single run.
$teamRange = $d.Range(130, 184)
$teamRange.Text = "Christian Ehringfeld, David Mock, Matthias Unterbusch"

# --- Change 2: "Protokollant:  Matthias Unterbusch" -> "Protokollant:  David Mock"
# Re-use the formatting (incl. <w:kern w:val="0"/>) already present on the
# "David Mock" text further up in the same paragraph, so the resulting run
# keeps the matching rPr instead of merging into the preceding "  " run.
$srcRange = $d.Range(152, 162)
$protokollantRange = $d.Range(202, 221)
$protokollantRange.FormattedText = $srcRange.FormattedText

Write-Output $d.Paragraphs.Item(8).Range.Text
Write-Output $d.Paragraphs.Item(11).Range.Text
